$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3703.3132
$ws.Range("J17").Value = 1309.4512
$ws.Range("L17").Value = 3928.3536
$ws.Range("N17").Value = -4264.3536

$ws.Range("H107").Value = 565.2
$ws.Range("I107").Value = 613.4737
$ws.Range("J107").Value = 412.33334
$ws.Range("K107").Value = 613.4737
$ws.Range("L107").Value = 412.33334
$ws.Range("M107").Value = 1306.5263
$ws.Range("N107").Value = -4252.33334

$ws.Range("H129").Value = 3013392.8
$ws.Range("I129").Value = 22728092
$ws.Range("J129").Value = 1424.8472
$ws.Range("K129").Value = 68184276
$ws.Range("L129").Value = 4274.5416
$ws.Range("M129").Value = -68179276
$ws.Range("N129").Value = -14274.5416

$ws.Range("H135").Value = 723.2
$ws.Range("I135").Value = 703.3333
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 6329.9997
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -3794.9997
$ws.Range("N135").Value = -15870

$ws.Range("H141").Value = 372079.06
$ws.Range("I141").Value = 1493.2703
$ws.Range("J141").Value = 2086038.4
$ws.Range("K141").Value = 4479.810899999999
$ws.Range("L141").Value = 6258115.199999999
$ws.Range("M141").Value = 700.1891000000005
$ws.Range("N141").Value = -6268475.199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1540.86
$ws.Range("I32").Value = 1230.1595
$ws.Range("J32").Value = 6408.5
$ws.Range("K32").Value = 1230.1595
$ws.Range("L32").Value = 6408.5
$ws.Range("M32").Value = -943.1595
$ws.Range("N32").Value = -6982.5

$ws.Range("H132").Value = 21278826
$ws.Range("I132").Value = 24391952
$ws.Range("J132").Value = 5799.8335
$ws.Range("K132").Value = 73175856
$ws.Range("L132").Value = 17399.5005
$ws.Range("M132").Value = -73173326
$ws.Range("N132").Value = -22459.5005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1964843.8
$ws.Range("I31").Value = 2860350
$ws.Range("J31").Value = 5923.9375
$ws.Range("K31").Value = 2860350
$ws.Range("L31").Value = 5923.9375
$ws.Range("M31").Value = -2860055
$ws.Range("N31").Value = -6513.9375

$ws.Range("H34").Value = 1964843.8
$ws.Range("I34").Value = 2860350
$ws.Range("J34").Value = 5923.9375
$ws.Range("K34").Value = 2860350
$ws.Range("L34").Value = 5923.9375
$ws.Range("M34").Value = -2860148
$ws.Range("N34").Value = -6327.9375

$ws.Range("H58").Value = 8774165
$ws.Range("I58").Value = 1472.8649
$ws.Range("J58").Value = 25003646
$ws.Range("K58").Value = 1472.8649
$ws.Range("L58").Value = 25003646
$ws.Range("M58").Value = -1269.8649
$ws.Range("N58").Value = -25004052

$ws.Range("H94").Value = 2344.375
$ws.Range("I94").Value = 2088.75
$ws.Range("K94").Value = 2088.75
$ws.Range("M94").Value = -1637.75

$ws.Range("H132").Value = 2077.2727
$ws.Range("I132").Value = 1509.0264
$ws.Range("J132").Value = 3347.4707
$ws.Range("K132").Value = 4527.0792
$ws.Range("L132").Value = 10042.4121
$ws.Range("M132").Value = -1997.0792
$ws.Range("N132").Value = -15102.4121

$ws.Range("H134").Value = 1147.6774
$ws.Range("I134").Value = 797.7222
$ws.Range("J134").Value = 1632.2307
$ws.Range("K134").Value = 2393.1666
$ws.Range("L134").Value = 4896.6921
$ws.Range("M134").Value = 141.8334
$ws.Range("N134").Value = -9966.6921

$ws.Range("H136").Value = 8774165
$ws.Range("I136").Value = 1472.8649
$ws.Range("J136").Value = 25003646
$ws.Range("K136").Value = 4418.5947
$ws.Range("L136").Value = 75010938
$ws.Range("M136").Value = -1868.5947
$ws.Range("N136").Value = -75016038

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 11209.091
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 12280
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 36840
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -37008

$ws.Range("H39").Value = 1650.5
$ws.Range("I39").Value = 450
$ws.Range("J39").Value = 1822
$ws.Range("K39").Value = 1350
$ws.Range("L39").Value = 5466
$ws.Range("M39").Value = -1056
$ws.Range("N39").Value = -6054

$ws.Range("H63").Value = 36000
$ws.Range("J63").Value = 36000
$ws.Range("L63").Value = 108000
$ws.Range("N63").Value = -109498

$ws.Range("H66").Value = 36000
$ws.Range("J66").Value = 36000
$ws.Range("L66").Value = 324000
$ws.Range("N66").Value = -331488

$ws.Range("H86").Value = 876.6667
$ws.Range("I86").Value = 630
$ws.Range("J86").Value = 1123.3334
$ws.Range("K86").Value = 1890
$ws.Range("L86").Value = 3370.0002
$ws.Range("M86").Value = -704
$ws.Range("N86").Value = -5742.0002

$ws.Range("H89").Value = 876.6667
$ws.Range("I89").Value = 630
$ws.Range("J89").Value = 1123.3334
$ws.Range("K89").Value = 5670
$ws.Range("L89").Value = 10110.0006
$ws.Range("M89").Value = 258
$ws.Range("N89").Value = -21966.0006

$ws.Range("H97").Value = 1687.5
$ws.Range("I97").Value = 1375
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 4125
$ws.Range("L97").Value = 6000
$ws.Range("M97").Value = -3629
$ws.Range("N97").Value = -6992

$ws.Range("H131").Value = 1207.9
$ws.Range("I131").Value = 2236.6667
$ws.Range("J131").Value = 1067.6136
$ws.Range("K131").Value = 6710.000100000001
$ws.Range("L131").Value = 3202.8408
$ws.Range("M131").Value = -1670.000100000001
$ws.Range("N131").Value = -13282.8408

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 645752
$ws.Range("J14").Value = 70005
$ws.Range("L14").Value = 70005
$ws.Range("N14").Value = -70349

$ws.Range("H16").Value = 2233.375
$ws.Range("I16").Value = 1123.8572
$ws.Range("J16").Value = 10000
$ws.Range("K16").Value = 1123.8572
$ws.Range("L16").Value = 10000
$ws.Range("M16").Value = -953.8571999999999
$ws.Range("N16").Value = -10340

$ws.Range("H22").Value = 83334530
$ws.Range("I22").Value = 250000290
$ws.Range("J22").Value = 1648.75
$ws.Range("K22").Value = 250000290
$ws.Range("L22").Value = 1648.75
$ws.Range("M22").Value = -249999995
$ws.Range("N22").Value = -2238.75

$ws.Range("H27").Value = 83334530
$ws.Range("I27").Value = 250000290
$ws.Range("J27").Value = 1648.75
$ws.Range("K27").Value = 250000290
$ws.Range("L27").Value = 1648.75
$ws.Range("M27").Value = -250000183
$ws.Range("N27").Value = -1862.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 24975
$ws.Range("J123").Value = 24975
$ws.Range("L123").Value = 24975
$ws.Range("N123").Value = -34775
